$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column D (was "Viso"), column E gets the "Viso" header that D used to have
$ws.Range("D1").Value = "DATE"
$ws.Range("E1").Value = "Viso"

# Move the per-row total formula from D to E, and put date text in D
$ws.Range("E2").Formula = "=B2*C2"
$ws.Range("E3").Formula = "=B3*C3"
$ws.Range("E4").Formula = "=B4*C4"
$ws.Range("E5").Formula = "=B5*C5"
$ws.Range("E6").Formula = "=B6*C6"
$ws.Range("E7").Formula = "=SUM(E2:E6)"

$ws.Range("D2").Value = "2021/10/24"
$ws.Range("D3").Value = "2024 - 12 - 10"
$ws.Range("D4").Value = "2021\10\24"
$ws.Range("D5").Value = "2021 - 10 - 24"
$ws.Range("D6").Value = "2021/10/24"
$ws.Range("D7").Value = "2021/10/24"
